$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$data = @(
    @(1.02,1.017860154984261,1.023792172426138,1.046603232189801,1.050278911179785,1.029272293840606,1.023071871476029,1.026622570117025,1.049368205708639,1.053033630876649,1.01179179074138),
    @(1.02,1.018653072087143,1.024345317934763,1.047630182861211,1.051314267519219,1.029371223387494,1.023501680761275,1.026983361458274,1.050206273090264,1.053880827045953,1.011934849430413),
    @(1.02,1.019166666266542,1.024703499657061,1.048295892362158,1.051985243310545,1.029434002432631,1.023779669073545,1.02721637177825,1.050749154597482,1.054429443945295,1.012027345339751),
    @(1.02,1.019382705088004,1.024854139680082,1.048576043728,1.052267567166968,1.029460098475218,1.023896503826951,1.027314221424714,1.050977523655827,1.05466018297034,1.012066212823809),
    @(1.02,1.019418986172532,1.02487943629618,1.048623099178752,1.052314984954315,1.029464462724589,1.023916119014441,1.027330644450144,1.051015876104969,1.054698930928176,1.012072737789932),
    @(1.02,1.01916955250325,1.024705512281154,1.048299634631701,1.051989014772045,1.029434352294121,1.023781230350549,1.027217679675423,1.050752205520562,1.054432526697999,1.012027864759416),
    @(1.02,1.018128015633248,1.023979055758099,1.046950045078071,1.050628600507247,1.029305982889699,1.023217153166532,1.026744592846055,1.049651310534224,1.053319856263462,1.011840152887739),
    @(1.02,1.016296775597048,1.022701010584787,1.04458116327676,1.048239325797189,1.029070348029807,1.022222253644609,1.025907590562884,1.047715999971992,1.051362487524463,1.011508842385202),
    @(1.02,1.015078785971291,1.021850470584718,1.043008203591826,1.046651890636119,1.028906953157379,1.021558433001564,1.025347397387447,1.046428948141566,1.05005985802854,1.011287631148837),
    @(1.02,1.014552075494914,1.021482552432591,1.042328600901717,1.045965812444464,1.028834713537922,1.021270873858945,1.025104322325013,1.04587240075148,1.049496360037029,1.011191769283273),
    @(1.02,1.014356536561561,1.0213459487802,1.042076392471236,1.045711167650365,1.028807657481865,1.021164044616794,1.025013958312579,1.04566578855496,1.049287135349576,1.011156150917967),
    @(1.02,1.014398475559456,1.021375248095908,1.04213048178009,1.045765780971978,1.028813471181433,1.021186960601979,1.025033345080206,1.045710102385994,1.049332010957267,1.011163791669456),
    @(1.02,1.014535910038956,1.021471259548355,1.042307748652428,1.045944759437897,1.028832481622496,1.021262043660237,1.025096854333056,1.04585531978614,1.049479063748487,1.011188825280123),
    @(1.02,1.014620601890136,1.021530423048647,1.042416998636491,1.046055059876761,1.02884416504484,1.02130830259486,1.025135974561732,1.045944808190635,1.049569678930006,1.01120424786846),
    @(1.02,1.015113756620234,1.021874896072725,1.043053338258629,1.046697450697222,1.028911716151467,1.021577514913179,1.025363518881437,1.046465900334925,1.050097267194315,1.01129399162009),
    @(1.02,1.015423284605628,1.022091075539272,1.043452899400785,1.047100752109539,1.028953691065352,1.021746353237554,1.025506116315367,1.046792970332784,1.050428357142655,1.011350265462016),
    @(1.02,1.015603893229879,1.022217205070733,1.043686101240156,1.047336115555086,1.028978030657955,1.021844822081497,1.025589241920099,1.046983817469494,1.050621529224824,1.011383081687388),
    @(1.02,1.015665487282306,1.022260217984538,1.043765641533973,1.047416389500118,1.028986305443968,1.021878395395023,1.025617577245163,1.047048903721259,1.05068740485663,1.011394269912281),
    @(1.02,1.015390068324789,1.022067877822529,1.043410015309085,1.047057468793306,1.028949202409941,1.021728239684715,1.025490822026919,1.046757871286791,1.05039282884356,1.011344228569108),
    @(1.02,1.014495436104931,1.021442984966467,1.042255541722418,1.045892049343196,1.028826889673778,1.021239934030529,1.025078154505087,1.045812553718999,1.049435758051652,1.011181453806105),
    @(1.02,1.013933552375296,1.021050423744525,1.041530987986155,1.045160433774957,1.028748696627663,1.020932819237784,1.024818260061258,1.045218856719822,1.048834493613147,1.011079047410402),
    @(1.02,1.014231359506762,1.021258495660013,1.041914963149152,1.04554816955055,1.028790270335764,1.021095635482526,1.024956075771396,1.045533523703546,1.049153189043383,1.011133340855495),
    @(1.02,1.015405077127998,1.022078359762089,1.043429392337526,1.04707702627603,1.028951231084078,1.02173642445582,1.02549773300928,1.046773730810161,1.050408882389197,1.011346956403932),
    @(1.02,1.016769701820554,1.023031160732023,1.04519247119905,1.048856061370973,1.029132379241967,1.022479561588844,1.026124367079934,1.048215772256944,1.051868116654308,1.011594555366235)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}

Write-Output "Updated $($data.Length) rows"
